$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2870
$ws.Range("J100").Value = 3916
$ws.Range("L100").Value = 3916
$ws.Range("N100").Value = -4998
$ws.Range("H112").Value = 5128.773
$ws.Range("J112").Value = 5128.773
$ws.Range("L112").Value = 15386.319
$ws.Range("N112").Value = -17602.319
$ws.Range("H138").Value = 4884.644
$ws.Range("I138").Value = 1476.9131
$ws.Range("J138").Value = 7061.8057
$ws.Range("K138").Value = 4430.7393
$ws.Range("L138").Value = 21185.4171
$ws.Range("M138").Value = 709.2606999999998
$ws.Range("N138").Value = -31465.4171

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 15250.444
$ws.Range("I31").Value = 6344.8335
$ws.Range("K31").Value = 6344.8335
$ws.Range("M31").Value = -6050.8335
$ws.Range("H45").Value = 1267.6111
$ws.Range("I45").Value = 758.5
$ws.Range("J45").Value = 3049.5
$ws.Range("K45").Value = 758.5
$ws.Range("L45").Value = 3049.5
$ws.Range("M45").Value = -381.5
$ws.Range("N45").Value = -3803.5
$ws.Range("H74").Value = 23257182
$ws.Range("I74").Value = 30303948
$ws.Range("K74").Value = 30303948
$ws.Range("M74").Value = -30303074
$ws.Range("H77").Value = 23257182
$ws.Range("I77").Value = 30303948
$ws.Range("K77").Value = 151519740
$ws.Range("M77").Value = -151515372
$ws.Range("H122").Value = 111114380

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 14248.4
$ws.Range("H107").Value = 1014.8
$ws.Range("I107").Value = 896.6667
$ws.Range("K107").Value = 896.6667
$ws.Range("M107").Value = 1023.3333
$ws.Range("H112").Value = 87666.664
$ws.Range("J112").Value = 84000
$ws.Range("L112").Value = 84000
$ws.Range("N112").Value = -86954

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1363.8462
$ws.Range("I58").Value = 1423.9
$ws.Range("J58").Value = 1163.6666
$ws.Range("K58").Value = 1423.9
$ws.Range("L58").Value = 1163.6666
$ws.Range("M58").Value = -1220.9
$ws.Range("N58").Value = -1569.6666
$ws.Range("H107").Value = 378.33334
$ws.Range("I107").Value = 324
$ws.Range("K107").Value = 324
$ws.Range("M107").Value = 1596
$ws.Range("H121").Value = 49333
$ws.Range("J121").Value = 36999.5
$ws.Range("L121").Value = 36999.5
$ws.Range("N121").Value = -39619.5
$ws.Range("H122").Value = 1951100.8
$ws.Range("I122").Value = 1594.1111
$ws.Range("K122").Value = 4782.3333
$ws.Range("M122").Value = -2332.3333
$ws.Range("H132").Value = 137098.94
$ws.Range("I132").Value = 204176.1
$ws.Range("J132").Value = 2944.6
$ws.Range("K132").Value = 612528.3
$ws.Range("L132").Value = 8833.799999999999
$ws.Range("M132").Value = -609998.3
$ws.Range("N132").Value = -13893.8
$ws.Range("H136").Value = 1363.8462
$ws.Range("I136").Value = 1423.9
$ws.Range("J136").Value = 1163.6666
$ws.Range("K136").Value = 4271.700000000001
$ws.Range("L136").Value = 3490.9998
$ws.Range("M136").Value = -1721.700000000001
$ws.Range("N136").Value = -8590.9998
$ws.Range("H140").Value = 93103.42
$ws.Range("I140").Value = 40709
$ws.Range("K140").Value = 40709
$ws.Range("M140").Value = -35529

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1273.1666
$ws.Range("J12").Value = 1333.5
$ws.Range("L12").Value = 4000.5
$ws.Range("N12").Value = -4346.5
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H61").Value = 178.38461
$ws.Range("I61").Value = 212.5
$ws.Range("K61").Value = 637.5
$ws.Range("M61").Value = -422.5
$ws.Range("H102").Value = 4398.6
$ws.Range("I102").Value = 2998
$ws.Range("K102").Value = 8994
$ws.Range("M102").Value = -6560
$ws.Range("H132").Value = 7412160.5
$ws.Range("I132").Value = 1355.8
$ws.Range("J132").Value = 16675666
$ws.Range("K132").Value = 12202.2
$ws.Range("L132").Value = 150080994
$ws.Range("M132").Value = -9672.199999999999
$ws.Range("N132").Value = -150086054

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1266.0714
$ws.Range("I22").Value = 949
$ws.Range("K22").Value = 949
$ws.Range("M22").Value = -654
$ws.Range("H27").Value = 1266.0714
$ws.Range("I27").Value = 949
$ws.Range("K27").Value = 949
$ws.Range("M27").Value = -842
$ws.Range("H82").Value = 1479.5555
$ws.Range("I82").Value = 1814.4
$ws.Range("J82").Value = 1061
$ws.Range("K82").Value = 1814.4
$ws.Range("L82").Value = 1061
$ws.Range("M82").Value = -1453.4
$ws.Range("N82").Value = -1783
$ws.Range("H85").Value = 1479.5555
$ws.Range("I85").Value = 1814.4
$ws.Range("J85").Value = 1061
$ws.Range("K85").Value = 1814.4
$ws.Range("L85").Value = 1061
$ws.Range("M85").Value = -566.4000000000001
$ws.Range("N85").Value = -3557
$ws.Range("H93").Value = 1519434
$ws.Range("I93").Value = 3887.8333
$ws.Range("K93").Value = 3887.8333
$ws.Range("M93").Value = -2639.8333
$ws.Range("H132").Value = 5739.1724
$ws.Range("I132").Value = 2655.4736
$ws.Range("K132").Value = 7966.4208
$ws.Range("M132").Value = -5436.4208

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 900
$ws.Range("I81").Value = 800
$ws.Range("J81").Value = 1000
$ws.Range("K81").Value = 1600
$ws.Range("L81").Value = 2000
$ws.Range("M81").Value = -539
$ws.Range("N81").Value = -4122
$ws.Range("H84").Value = 900
$ws.Range("I84").Value = 800
$ws.Range("J84").Value = 1000
$ws.Range("K84").Value = 8000
$ws.Range("L84").Value = 10000
$ws.Range("M84").Value = -2696
$ws.Range("N84").Value = -20608
$ws.Range("H96").Value = 4750.4287
$ws.Range("J96").Value = 5377
$ws.Range("L96").Value = 5377
$ws.Range("N96").Value = -8123
$ws.Range("H107").Value = 2096.3333
$ws.Range("I107").Value = 777.8570999999999
$ws.Range("J107").Value = 3250
$ws.Range("K107").Value = 2333.5713
$ws.Range("L107").Value = 9750
$ws.Range("M107").Value = -413.5712999999996
$ws.Range("N107").Value = -13590
